$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 11.20127533333333
$ws.Range("H2").Value = 33.603826
$ws.Range("I2").Value = 0.1186573945858706
$ws.Range("J2").Value = 0.1186573945858706
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 13.441269
$ws.Range("N2").Value = 40.323807
$ws.Range("O2").Value = 0.08973082133481231
$ws.Range("P2").Value = 0.08973082133481232
$ws.Range("Q2").Value = 150.559354898398
$ws.Range("R2").Value = 1355.034194085582
$ws.Range("S2").Value = 0.01064722547363908
$ws.Range("T2").Value = 0.01064722547363908

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 11.20127533333333
$ws.Range("H3").Value = 33.603826
$ws.Range("I3").Value = 0.1186573945858706
$ws.Range("J3").Value = 0.1186573945858706
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 54.711535
$ws.Range("N3").Value = 164.134605
$ws.Range("O3").Value = 0.3652416280068742
$ws.Range("P3").Value = 0.3652416280068742
$ws.Range("Q3").Value = 612.8389674443033
$ws.Range("R3").Value = 5515.55070699873
$ws.Range("S3").Value = 0.04333861997359744
$ws.Range("T3").Value = 0.04333861997359745

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 11.20127533333333
$ws.Range("H4").Value = 33.603826
$ws.Range("I4").Value = 0.1186573945858706
$ws.Range("J4").Value = 0.1186573945858706
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 63.67711
$ws.Range("N4").Value = 191.03133
$ws.Range("O4").Value = 0.4250937452800914
$ws.Range("P4").Value = 0.4250937452800915
$ws.Range("Q4").Value = 713.2648415409533
$ws.Range("R4").Value = 6419.383573868579
$ws.Range("S4").Value = 0.05044051626968539
$ws.Range("T4").Value = 0.0504405162696854

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 11.20127533333333
$ws.Range("H5").Value = 33.603826
$ws.Range("I5").Value = 0.1186573945858706
$ws.Range("J5").Value = 0.1186573945858706
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.96553866666667
$ws.Range("N5").Value = 53.896616
$ws.Range("O5").Value = 0.119933805378222
$ws.Range("P5").Value = 0.119933805378222
$ws.Range("Q5").Value = 201.2369451169796
$ws.Range("R5").Value = 1811.132506052816
$ws.Range("S5").Value = 0.0142310328689487
$ws.Range("T5").Value = 0.0142310328689487

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 47.94465366666667
$ws.Range("H6").Value = 143.833961
$ws.Range("I6").Value = 0.5078874966566524
$ws.Range("J6").Value = 0.5078874966566524
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 13.441269
$ws.Range("N6").Value = 40.323807
$ws.Range("O6").Value = 0.08973082133481231
$ws.Range("P6").Value = 0.08973082133481232
$ws.Range("Q6").Value = 644.4369870455031
$ws.Range("R6").Value = 5799.932883409528
$ws.Range("S6").Value = 0.04557316222068316
$ws.Range("T6").Value = 0.04557316222068317

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 47.94465366666667
$ws.Range("H7").Value = 143.833961
$ws.Range("I7").Value = 0.5078874966566524
$ws.Range("J7").Value = 0.5078874966566524
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 54.711535
$ws.Range("N7").Value = 164.134605
$ws.Range("O7").Value = 0.3652416280068742
$ws.Range("P7").Value = 0.3652416280068742
$ws.Range("Q7").Value = 2623.125597146712
$ws.Range("R7").Value = 23608.13037432041
$ws.Range("S7").Value = 0.1855016561232116
$ws.Range("T7").Value = 0.1855016561232116

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 47.94465366666667
$ws.Range("H8").Value = 143.833961
$ws.Range("I8").Value = 0.5078874966566524
$ws.Range("J8").Value = 0.5078874966566524
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 63.67711
$ws.Range("N8").Value = 191.03133
$ws.Range("O8").Value = 0.4250937452800914
$ws.Range("P8").Value = 0.4250937452800915
$ws.Range("Q8").Value = 3052.976985444237
$ws.Range("R8").Value = 27476.79286899813
$ws.Range("S8").Value = 0.2158997981347063
$ws.Range("T8").Value = 0.2158997981347063

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 47.94465366666667
$ws.Range("H9").Value = 143.833961
$ws.Range("I9").Value = 0.5078874966566524
$ws.Range("J9").Value = 0.5078874966566524
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 17.96553866666667
$ws.Range("N9").Value = 53.896616
$ws.Range("O9").Value = 0.119933805378222
$ws.Range("P9").Value = 0.119933805378222
$ws.Range("Q9").Value = 861.351529308442
$ws.Range("R9").Value = 7752.163763775977
$ws.Range("S9").Value = 0.06091288017805133
$ws.Range("T9").Value = 0.06091288017805134

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 25.59984766666667
$ws.Range("H10").Value = 76.799543
$ws.Range("I10").Value = 0.2711844085184091
$ws.Range("J10").Value = 0.2711844085184091
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 13.441269
$ws.Range("N10").Value = 40.323807
$ws.Range("O10").Value = 0.08973082133481231
$ws.Range("P10").Value = 0.08973082133481232
$ws.Range("Q10").Value = 344.094438846689
$ws.Range("R10").Value = 3096.849949620201
$ws.Range("S10").Value = 0.02433359970955212
$ws.Range("T10").Value = 0.02433359970955213

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 25.59984766666667
$ws.Range("H11").Value = 76.799543
$ws.Range("I11").Value = 0.2711844085184091
$ws.Range("J11").Value = 0.2711844085184091
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 54.711535
$ws.Range("N11").Value = 164.134605
$ws.Range("O11").Value = 0.3652416280068742
$ws.Range("P11").Value = 0.3652416280068742
$ws.Range("Q11").Value = 1400.606961609501
$ws.Range("R11").Value = 12605.46265448551
$ws.Range("S11").Value = 0.09904783485734499
$ws.Range("T11").Value = 0.099047834857345

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 25.59984766666667
$ws.Range("H12").Value = 76.799543
$ws.Range("I12").Value = 0.2711844085184091
$ws.Range("J12").Value = 0.2711844085184091
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 63.67711
$ws.Range("N12").Value = 191.03133
$ws.Range("O12").Value = 0.4250937452800914
$ws.Range("P12").Value = 0.4250937452800915
$ws.Range("Q12").Value = 1630.124315853576
$ws.Range("R12").Value = 14671.11884268219
$ws.Range("S12").Value = 0.1152787958786569
$ws.Range("T12").Value = 0.1152787958786569

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 25.59984766666667
$ws.Range("H13").Value = 76.799543
$ws.Range("I13").Value = 0.2711844085184091
$ws.Range("J13").Value = 0.2711844085184091
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 17.96553866666667
$ws.Range("N13").Value = 53.896616
$ws.Range("O13").Value = 0.119933805378222
$ws.Range("P13").Value = 0.119933805378222
$ws.Range("Q13").Value = 459.9150531162765
$ws.Range("R13").Value = 4139.235478046488
$ws.Range("S13").Value = 0.03252417807285513
$ws.Range("T13").Value = 0.03252417807285513

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 9.654369000000001
$ws.Range("H14").Value = 28.963107
$ws.Range("I14").Value = 0.1022707002390678
$ws.Range("J14").Value = 0.1022707002390678
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 13.441269
$ws.Range("N14").Value = 40.323807
$ws.Range("O14").Value = 0.08973082133481231
$ws.Range("P14").Value = 0.08973082133481232
$ws.Range("Q14").Value = 129.766970754261
$ws.Range("R14").Value = 1167.902736788349
$ws.Range("S14").Value = 0.009176833930937938
$ws.Range("T14").Value = 0.009176833930937938

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 9.654369000000001
$ws.Range("H15").Value = 28.963107
$ws.Range("I15").Value = 0.1022707002390678
$ws.Range("J15").Value = 0.1022707002390678
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 54.711535
$ws.Range("N15").Value = 164.134605
$ws.Range("O15").Value = 0.3652416280068742
$ws.Range("P15").Value = 0.3652416280068742
$ws.Range("Q15").Value = 528.205347446415
$ws.Range("R15").Value = 4753.848127017735
$ws.Range("S15").Value = 0.03735351705272013
$ws.Range("T15").Value = 0.03735351705272013

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 9.654369000000001
$ws.Range("H16").Value = 28.963107
$ws.Range("I16").Value = 0.1022707002390678
$ws.Range("J16").Value = 0.1022707002390678
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 63.67711
$ws.Range("N16").Value = 191.03133
$ws.Range("O16").Value = 0.4250937452800914
$ws.Range("P16").Value = 0.4250937452800915
$ws.Range("Q16").Value = 614.76231679359
$ws.Range("R16").Value = 5532.86085114231
$ws.Range("S16").Value = 0.04347463499704286
$ws.Range("T16").Value = 0.04347463499704286

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 9.654369000000001
$ws.Range("H17").Value = 28.963107
$ws.Range("I17").Value = 0.1022707002390678
$ws.Range("J17").Value = 0.1022707002390678
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 17.96553866666667
$ws.Range("N17").Value = 53.896616
$ws.Range("O17").Value = 0.119933805378222
$ws.Range("P17").Value = 0.119933805378222
$ws.Range("Q17").Value = 173.445939571768
$ws.Range("R17").Value = 1561.013456145912
$ws.Range("S17").Value = 0.01226571425836684
$ws.Range("T17").Value = 0.01226571425836684
